$d = $word.ActiveDocument

# Add a new paragraph style "Compact List" (styleId "CompactList"), mirroring
# the existing "Compact" style: based on Body Text, quick-style, with
# before/after paragraph spacing of 36 twips (1.8pt).
$style = $d.Styles.Add("Compact List", 1)
$style.BaseStyle = "BodyText"
$style.QuickStyle = $true
$style.ParagraphFormat.SpaceBefore = 1.8
$style.ParagraphFormat.SpaceAfter = 1.8
